$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45202 -> 45203, i.e. 2023-10-03 -> 2023-10-04) for every data row (2..262).
$ws.Range("C2:C262").Value = 45203
